$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") "30.208.32"
Set-TextValue $ws.Range("E2") "  -0.64%  "

# Row 3
Set-TextValue $ws.Range("D3") "1.909.01"
Set-TextValue $ws.Range("E3") "  -1.52%  "

# Row 4
Set-TextValue $ws.Range("D4") "1.002"
Set-TextValue $ws.Range("E4") "  +0.06%  "

# Row 5
Set-TextValue $ws.Range("D5") "0.7332"
Set-TextValue $ws.Range("E5") "  -4.27%  "

# Row 6
Set-TextValue $ws.Range("D6") "243.73"
Set-TextValue $ws.Range("E6") "  -1.61%  "

# Row 7
Set-TextValue $ws.Range("D7") "1.002"
Set-TextValue $ws.Range("E7") "  +0.18%  "

# Row 8
Set-TextValue $ws.Range("D8") "0.3134"
Set-TextValue $ws.Range("E8") "  -2.03%  "

# Row 9
Set-TextValue $ws.Range("E9") "  -3.98%  "

# Row 10
Set-TextValue $ws.Range("D10") "0.06918"
Set-TextValue $ws.Range("E10") "  -2.63%  "

# Row 11
Set-TextValue $ws.Range("D11") "0.7781"
Set-TextValue $ws.Range("E11") "  -0.79%  "

# Row 12
Set-TextValue $ws.Range("D12") "0.07978"
Set-TextValue $ws.Range("E12") "  -0.47%  "

# Row 13
Set-TextValue $ws.Range("D13") "1.900.41"
Set-TextValue $ws.Range("E13") "  -1.93%  "

# Row 14
Set-TextValue $ws.Range("D14") "5.254"
Set-TextValue $ws.Range("E14") "  -2.25%  "

# Row 15
Set-TextValue $ws.Range("D15") "91.48"
Set-TextValue $ws.Range("E15") "  -3.79%  "

# Row 16
Set-TextValue $ws.Range("D16") "30.141.20"
Set-TextValue $ws.Range("E16") "  -0.85%  "

# Row 17
Set-TextValue $ws.Range("E17") "  -2.23%  "

# Row 18
Set-TextValue $ws.Range("D18") "5.852"
Set-TextValue $ws.Range("E18") "  +0.05%  "

# Row 19
Set-TextValue $ws.Range("D19") "240.34"
Set-TextValue $ws.Range("E19") "  -6.87%  "

# Row 20
Set-TextValue $ws.Range("D20") "0.000007800"
Set-TextValue $ws.Range("E20") "  -2.59%  "

# Row 21
Set-TextValue $ws.Range("E21") "  +0.19%  "

# Row 22
Set-TextValue $ws.Range("D22") "2.134.89"
Set-TextValue $ws.Range("E22") "  -2.49%  "

# Row 23
Set-TextValue $ws.Range("D23") "1.002"
Set-TextValue $ws.Range("E23") "  -0.04%  "

# Row 24
Set-TextValue $ws.Range("D24") "6.805"
Set-TextValue $ws.Range("E24") "  +0.62%  "

# Row 25
Set-TextValue $ws.Range("D25") "9.391"
Set-TextValue $ws.Range("E25") "  -2.36%  "

# Row 26
Set-TextValue $ws.Range("D26") "165.76"
Set-TextValue $ws.Range("E26") "  +0.71%  "

# Row 27
Set-TextValue $ws.Range("D27") "19.11"
Set-TextValue $ws.Range("E27") "  -0.23%  "

# Row 28
Set-TextValue $ws.Range("D28") "0.1272"
Set-TextValue $ws.Range("E28") "  -4.37%  "

# Row 29
Set-TextValue $ws.Range("E29") "  -8.95%  "

# Row 30
Set-TextValue $ws.Range("B30") "Toncoin"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
Set-TextValue $ws.Range("D30") "1.347"
Set-TextValue $ws.Range("E30") "  -1.19%  "

# Row 31
Set-TextValue $ws.Range("B31") "PancakeSwap"
Set-TextValue $ws.Range("C31") "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
Set-TextValue $ws.Range("D31") "1.548"
Set-TextValue $ws.Range("E31") "  +1.49%  "

# Row 32
Set-TextValue $ws.Range("D32") "4.305"
Set-TextValue $ws.Range("E32") "  -3.04%  "

# Row 33
Set-TextValue $ws.Range("D33") "4.086"
Set-TextValue $ws.Range("E33") "  -1.47%  "

# Row 34
Set-TextValue $ws.Range("D34") "0.05163"
Set-TextValue $ws.Range("E34") "  -0.72%  "

# Row 35
Set-TextValue $ws.Range("D35") "1.284"
Set-TextValue $ws.Range("E35") "  +0.24%  "

# Row 36
Set-TextValue $ws.Range("D36") "0.7423"
Set-TextValue $ws.Range("E36") "  -1.00%  "

# Row 37
Set-TextValue $ws.Range("D37") "2.753"
Set-TextValue $ws.Range("E37") "  -0.92%  "

# Row 38
Set-TextValue $ws.Range("D38") "0.01941"
Set-TextValue $ws.Range("E38") "  -1.45%  "

# Row 39
Set-TextValue $ws.Range("E39") "  -0.53%  "

# Row 40
Set-TextValue $ws.Range("E40") "  -1.81%  "

# Row 41
Set-TextValue $ws.Range("D41") "74.76"
Set-TextValue $ws.Range("E41") "  -4.75%  "

# Row 42
Set-TextValue $ws.Range("D42") "0.4440"
Set-TextValue $ws.Range("E42") "  -2.02%  "

# Row 43
Set-TextValue $ws.Range("D43") "1.931"
Set-TextValue $ws.Range("E43") "  -2.14%  "

# Row 44
Set-TextValue $ws.Range("E44") "  +0.01%  "

# Row 45
Set-TextValue $ws.Range("D45") "0.8343"
Set-TextValue $ws.Range("E45") "  -0.03%  "

# Row 46
Set-TextValue $ws.Range("D46") "101.04"
Set-TextValue $ws.Range("E46") "  -0.32%  "

# Row 47
Set-TextValue $ws.Range("D47") "7.588"
Set-TextValue $ws.Range("E47") "  +0.73%  "

# Row 48
Set-TextValue $ws.Range("D48") "9.728"
Set-TextValue $ws.Range("E48") "  -0.70%  "

# Row 49
Set-TextValue $ws.Range("D49") "37.52"
Set-TextValue $ws.Range("E49") "  +0.36%  "

# Row 50
Set-TextValue $ws.Range("B50") "Maker"
Set-TextValue $ws.Range("C50") "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
Set-TextValue $ws.Range("D50") "945.58"
Set-TextValue $ws.Range("E50") "  -4.08%  "

# Row 51
Set-TextValue $ws.Range("B51") "Algorand"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws.Range("D51") "0.1185"
Set-TextValue $ws.Range("E51") "  +0.54%  "
